# Security Task Team Information.xlsx
# Commit: "add team name and repo"
#
# The sheet holds one header row (A1:G1 = team-member-slot labels, A1 = repo
# link label) and one data row (A2:G2). This edit:
#   - replaces the repo link + all team member names in row 2 with a new
#     team's data (A2's cell keeps its existing hyperlink relationship -
#     only the *displayed* text/shared-string changes, matching the diff
#     which leaves xl/worksheets/_rels/sheet1.xml.rels untouched)
#   - adds a 7th member column (H): header "اسم الفرد السابع" + a new name
#   - gives the new/changed header cell (G1, which previously had a
#     slightly different "odd one out" font) and the new H1 the same
#     look as the rest of the header row, and gives H2 the same look as
#     the rest of the data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new repo link + new team member names -------------------------
# Setting .Value on A2 (which already carries a hyperlink) updates the
# displayed text in place while leaving the hyperlink target untouched.
$ws.Range("A2").Value = "https://github.com/Momen-magdy/websiteforsianna"
$ws.Range("B2").Value = "مؤمن مجدى عبدالعزيز"
$ws.Range("C2").Value = "احمد بلال عبدالمجيد"
$ws.Range("D2").Value = "احمد اسماعيل دياب"
$ws.Range("E2").Value = "محمود محمد فؤاد"
$ws.Range("F2").Value = "محمود محمد ابراهيم"
$ws.Range("G2").Value = "عمر محمد محمد "

# --- New 7th member column --------------------------------------------------
$ws.Range("H1").Value = "اسم الفرد السابع"
$ws.Range("H2").Value = "احمد يحى محمد "

# --- Normalize formatting ---------------------------------------------------
# G1 used to look different from the other header cells; make it (and the
# brand-new H1) match the rest of the blue header row.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null

# H2 should look like the rest of the (centered) data row.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Put the visible selection on the newly-added H1 cell.
$ws.Range("H1").Select() | Out-Null
